$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 112/113, pushing the existing rows 112..139 down to 114..141.
$ws.Rows.Item(112).Insert()
$ws.Rows.Item(112).Insert()

# Row 112: new Pimiento "Zafiro rojo" record
$ws.Range("A112").Value = 7
$ws.Range("B112").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C112").Value = "Ñuble"
$ws.Range("D112").Value = 44463
$ws.Range("E112").Value = 16
$ws.Range("F112").Value = 100112002
$ws.Range("G112").Value = "Pimiento"
$ws.Range("H112").Value = "Zafiro rojo"
$ws.Range("I112").Value = "Primera"
$ws.Range("J112").Value = 100
$ws.Range("K112").Value = 3800
$ws.Range("L112").Value = 39000
$ws.Range("M112").Value = 21400
$ws.Range("N112").Value = '$/caja 15 kilos'
$ws.Range("O112").Value = "Región de Arica y Parinacota"
$ws.Range("P112").Value = 1427
$ws.Range("Q112").Value = 15
$ws.Range("R112").Value = "Hortaliza"

# Row 113: new Pimiento "Zafiro verde" record
$ws.Range("A113").Value = 7
$ws.Range("B113").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C113").Value = "Ñuble"
$ws.Range("D113").Value = 44463
$ws.Range("E113").Value = 16
$ws.Range("F113").Value = 100112002
$ws.Range("G113").Value = "Pimiento"
$ws.Range("H113").Value = "Zafiro verde"
$ws.Range("I113").Value = "Primera"
$ws.Range("J113").Value = 160
$ws.Range("K113").Value = 29000
$ws.Range("L113").Value = 30000
$ws.Range("M113").Value = 29500
$ws.Range("N113").Value = '$/caja 15 kilos'
$ws.Range("O113").Value = "Región de Arica y Parinacota"
$ws.Range("P113").Value = 1967
$ws.Range("Q113").Value = 15
$ws.Range("R113").Value = "Hortaliza"
